# Updated IPS AIP hipo turnover
# Applies the per-location "Internal Fill Rate" / turnover corrections
# described in the commit: several monthly/quarterly/FY roll-up values
# on a handful of location sheets are refreshed, and a few cells that had
# a stray 0/avg value are cleared back to blank.

$wb = $excel.ActiveWorkbook

function Set-Cell {
    param($sheetName, $cellRef, $value)
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $value
}

function Clear-Cell {
    param($sheetName, $cellRef)
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).ClearContents()
}

# --- Manila Philippines (Internal Fill Rate, Commit/Forecast row 5) ---
Clear-Cell "Manila Philippines" "M5"
Set-Cell   "Manila Philippines" "N5" 1

# --- Milwaukee Pmc Hq Wisconsin ---
Set-Cell "Milwaukee Pmc Hq Wisconsin" "E2" 0.1754
Set-Cell "Milwaukee Pmc Hq Wisconsin" "E3" 0.1754
Set-Cell "Milwaukee Pmc Hq Wisconsin" "E4" 0.1754

Set-Cell "Milwaukee Pmc Hq Wisconsin" "M4" 0
Set-Cell "Milwaukee Pmc Hq Wisconsin" "N4" 0
Set-Cell "Milwaukee Pmc Hq Wisconsin" "O4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "P4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "Q4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "R4" 0.0877
Set-Cell "Milwaukee Pmc Hq Wisconsin" "S4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "T4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "U4" 0.0292333333333333
Set-Cell "Milwaukee Pmc Hq Wisconsin" "V4" 0.0877
Set-Cell "Milwaukee Pmc Hq Wisconsin" "W4" 0.3508

Clear-Cell "Milwaukee Pmc Hq Wisconsin" "M7"
Set-Cell   "Milwaukee Pmc Hq Wisconsin" "N7" 0

# --- Monterrey Rbm Mexico ---
Set-Cell "Monterrey Rbm Mexico" "E2" 0.1852
Set-Cell "Monterrey Rbm Mexico" "E3" 0.1852
Set-Cell "Monterrey Rbm Mexico" "E4" 0.1852

Set-Cell "Monterrey Rbm Mexico" "M4" 0
Set-Cell "Monterrey Rbm Mexico" "N4" 0.2062
Set-Cell "Monterrey Rbm Mexico" "O4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "P4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "Q4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "R4" 0.0926
Set-Cell "Monterrey Rbm Mexico" "S4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "T4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "U4" 0.0308666666666667
Set-Cell "Monterrey Rbm Mexico" "V4" 0.0926
Set-Cell "Monterrey Rbm Mexico" "W4" 0.3704

Set-Cell "Monterrey Rbm Mexico" "E7" 0.3043
Set-Cell "Monterrey Rbm Mexico" "E8" 0.3043

Set-Cell "Monterrey Rbm Mexico" "E9" 0.3043
Set-Cell "Monterrey Rbm Mexico" "G9" 0.0435
Set-Cell "Monterrey Rbm Mexico" "H9" 0.1364
Set-Cell "Monterrey Rbm Mexico" "J9" 0.1762
Set-Cell "Monterrey Rbm Mexico" "K9" 0.0417
Set-Cell "Monterrey Rbm Mexico" "L9" 0.0435
Set-Cell "Monterrey Rbm Mexico" "M9" 0.0435
Set-Cell "Monterrey Rbm Mexico" "N9" 0.1288
Set-Cell "Monterrey Rbm Mexico" "O9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "P9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "Q9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "R9" 0.15215
Set-Cell "Monterrey Rbm Mexico" "S9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "T9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "U9" 0.0507166666666667
Set-Cell "Monterrey Rbm Mexico" "V9" 0.15215
Set-Cell "Monterrey Rbm Mexico" "W9" 0.6086

# --- Rosemont Illinois ---
Clear-Cell "Rosemont Illinois" "M7"

# --- Tipp City Ohio ---
Set-Cell "Tipp City Ohio" "E2" 0.6667
Set-Cell "Tipp City Ohio" "E3" 0.6667
Set-Cell "Tipp City Ohio" "E4" 0.6667

Set-Cell "Tipp City Ohio" "M4" 0
Set-Cell "Tipp City Ohio" "N4" 1
Set-Cell "Tipp City Ohio" "O4" 0.111116666666667
Set-Cell "Tipp City Ohio" "P4" 0.111116666666667
Set-Cell "Tipp City Ohio" "Q4" 0.111116666666667
Set-Cell "Tipp City Ohio" "R4" 0.33335
Set-Cell "Tipp City Ohio" "S4" 0.111116666666667
Set-Cell "Tipp City Ohio" "T4" 0.111116666666667
Set-Cell "Tipp City Ohio" "U4" 0.111116666666667
Set-Cell "Tipp City Ohio" "V4" 0.33335
Set-Cell "Tipp City Ohio" "W4" 1.3334

# --- Guadalajara Mexico ---
Clear-Cell "Guadalajara Mexico" "M4"

# --- Faridabad India ---
Clear-Cell "Faridabad India" "M4"
Clear-Cell "Faridabad India" "N4"
